$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule" ----
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = -90.22366275
$wsSchedule.Range("F2").Value = -1.989057820767196
$wsSchedule.Range("E3").Value = 370.293807
$wsSchedule.Range("F3").Value = 24.49033115079365

# ---- Sheet "Detailed" ----
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B14").Value = 57.06003
$wsDetailed.Range("B16").Value = 0.7
$wsDetailed.Range("B17").Value = 0.51
$wsDetailed.Range("C17").Value = "historical"
$wsDetailed.Range("B18").Value = -4.34103
$wsDetailed.Range("C18").Value = "historical"
$wsDetailed.Range("B19").Value = 11.098
$wsDetailed.Range("B20").Value = 0.51
$wsDetailed.Range("B21").Value = -4.61115
$wsDetailed.Range("B22").Value = -0.90326
$wsDetailed.Range("B23").Value = -4.83666
$wsDetailed.Range("B24").Value = -5.51
$wsDetailed.Range("B25").Value = -4.5561
$wsDetailed.Range("B26").Value = -5.19013
$wsDetailed.Range("B27").Value = -5.01
$wsDetailed.Range("B28").Value = -5.51011
$wsDetailed.Range("B29").Value = -5.51011
$wsDetailed.Range("B30").Value = -7.01
$wsDetailed.Range("B31").Value = -12.08785
$wsDetailed.Range("B32").Value = -14.70709
$wsDetailed.Range("B33").Value = -14
$wsDetailed.Range("B34").Value = -6.85852
$wsDetailed.Range("B35").Value = -7.29618
$wsDetailed.Range("B37").Value = 0.66345
$wsDetailed.Range("B38").Value = 3.98398
$wsDetailed.Range("B39").Value = 3.93567
$wsDetailed.Range("B40").Value = 36.2
$wsDetailed.Range("B41").Value = 53.77198
$wsDetailed.Range("B43").Value = 56.04515
$wsDetailed.Range("B45").Value = 45.37537
$wsDetailed.Range("B46").Value = 55.9078
